$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (Test Case 5): "Expected Behavior" (col C) updated ---
$ws.Range("C6").Value = "Users cannot view any projects"

# --- Row 7 (Test Case 6): "Expected Behavior" (col C) updated ---
$ws.Range("C7").Value = "Users can only view for 2 room projects"

# --- Row 8 (Test Case 7): "Expected Behavior" (col C) updated ---
$ws.Range("C8").Value = "Users can view both 2 and 3 room projects"

# --- Row 9 becomes new Test Case 8 (was a mostly-empty placeholder row) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Project Application`nUser aged >= 35 and Single: `nOnly eligible for 2 room"
$ws.Range("C9").Value = "Users can only apply for 2 room projects"
$ws.Range("E9").Value = "1. Login as Applicant: Tom`n2. Enter NRIC: T1213151J`n3. Enter Pasword: password1`n4. Enter 1 for Projects`n5. Enter 1 to view list of projects`n6. Enter 2 to apply for for project`n7. Enter project name: Toa Payoh Glades"
$ws.Range("F9").Value = "Successfully applied for Toa Payoh Glades."
$ws.Range("G9").Value = "Successfully applied for Toa Payoh Glades."
$ws.Range("H9").Value = "Pass"
$ws.Rows.Item(9).RowHeight = 108.5

# --- Row 10 becomes new Test Case 9 (was a mostly-empty placeholder row) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Project Application`nUser aged > = 21 and Married: `nOnly eligible for 2 room"
$ws.Range("C10").Value = "Users can apply for both 2 and 3 room projects"
$ws.Range("E10").Value = "1. Login as Applicant: Ray`n2. Enter NRIC: T1122334K`n3. Enter Pasword: password1`n4. Enter 1 for Projects`n5. Enter 1 to view list of projects`n6. Enter 2 to apply for for project`n7. Enter 2 to apply for 3 room`n8. Enter project name: Toa Payoh Glades"
$ws.Range("F10").Value = "Successfully applied for Toa Payoh Glades."
$ws.Range("G10").Value = "Successfully applied for Toa Payoh Glades."
$ws.Range("H10").Value = "Pass"
$ws.Rows.Item(10).RowHeight = 124

# --- Update view state to match where the author left off scrolling/selecting ---
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("F12").Select()
